$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Update the 2014 season row (previously pointed at the 2013 file/label)
$ws.Range("A2").Value = "C:\Users\zaka\Desktop\MOTOGP\Excels\data\2014.xlsx"
$ws.Range("B2").Value = "2014"

$ws.Range("A2").Select()
